# paises.xlsx update: refresh COVID case counters and fix country-name
# ordering (4 rank swaps caused by the data refresh) plus the "last updated" timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 4: Estados Unidos - refreshed counters
$ws.Range("B4").Value = 2311824
$ws.Range("C4").Value = 14634
$ws.Range("D4").Value = 957204
$ws.Range("E4").Value = 1233018
$ws.Range("G4").Value = 195
$ws.Range("H4").Value = 121602

# Row 7: India - refreshed counters
$ws.Range("B7").Value = 402287
$ws.Range("C7").Value = 6475
$ws.Range("D7").Value = 218268
$ws.Range("E7").Value = 170941
$ws.Range("G7").Value = 108
$ws.Range("H7").Value = 13078

# Row 9: España - refreshed counters
$ws.Range("B9").Value = 293018
$ws.Range("C9").Value = 363
$ws.Range("G9").Value = 7
$ws.Range("H9").Value = 28322

# Row 11: Italia - refreshed counters
$ws.Range("B11").Value = 238275
$ws.Range("C11").Value = 264
$ws.Range("D11").Value = 182453
$ws.Range("E11").Value = 21212
$ws.Range("G11").Value = 49
$ws.Range("H11").Value = 34610

# Row 12: Chile - refreshed counters
$ws.Range("B12").Value = 236748
$ws.Range("C12").Value = 5355
$ws.Range("E12").Value = 40962
$ws.Range("G12").Value = 202
$ws.Range("H12").Value = 4295

# Row 14: Alemania - refreshed counters
$ws.Range("B14").Value = 190836
$ws.Range("C14").Value = 176
$ws.Range("E14").Value = 7476

# Row 35: Argentina - refreshed counters
$ws.Range("D35").Value = 12206
$ws.Range("E35").Value = 26384
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 980

# Row 42: Oman -> Irak (rank swap)
$ws.Range("A42").Value = "Irak"
$ws.Range("B42").Value = 29222
$ws.Range("C42").Value = 1870
$ws.Range("D42").Value = 13211
$ws.Range("E42").Value = 14998
$ws.Range("G42").Value = 88
$ws.Range("H42").Value = 1013

# Row 43: Afganistan -> Oman (rank swap)
$ws.Range("A43").Value = "Oman"
$ws.Range("B43").Value = 28566
$ws.Range("C43").Value = 896
$ws.Range("D43").Value = 14780
$ws.Range("E43").Value = 13658
$ws.Range("G43").Value = 3
$ws.Range("H43").Value = 128

# Row 44: Irak -> Afganistan (rank swap)
$ws.Range("A44").Value = "Afganistan"
$ws.Range("B44").Value = 28424
$ws.Range("C44").Value = 546
$ws.Range("D44").Value = 8292
$ws.Range("E44").Value = 19563
$ws.Range("G44").Value = 21
$ws.Range("H44").Value = 569

# Row 56: Kazajistan - refreshed counters
$ws.Range("D56").Value = 10671
$ws.Range("E56").Value = 5995

# Row 64: Argelia - refreshed counters
$ws.Range("B64").Value = 11631
$ws.Range("C64").Value = 127
$ws.Range("D64").Value = 8324
$ws.Range("E64").Value = 2470
$ws.Range("G64").Value = 12
$ws.Range("H64").Value = 837

# Row 67: Chequia - refreshed counters
$ws.Range("B67").Value = 10429
$ws.Range("C67").Value = 23
$ws.Range("D67").Value = 7477
$ws.Range("E67").Value = 2616
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 336

# Row 72: Sudan - refreshed counters
$ws.Range("B72").Value = 8416
$ws.Range("C72").Value = 100
$ws.Range("D72").Value = 3204
$ws.Range("E72").Value = 4699
$ws.Range("G72").Value = 7
$ws.Range("H72").Value = 513

# Row 76: Uzbekistan - refreshed counters
$ws.Range("D76").Value = 4290
$ws.Range("E76").Value = 1810

# Row 86: Etiopia - refreshed counters
$ws.Range("D86").Value = 1122
$ws.Range("E86").Value = 3275

# Row 88: Luxemburgo - refreshed counters
$ws.Range("B88").Value = 4105
$ws.Range("C88").Value = 6
$ws.Range("D88").Value = 3951
$ws.Range("E88").Value = 44

# Row 99: Mayotte - refreshed counters
$ws.Range("B99").Value = 2404
$ws.Range("C99").Value = 10
$ws.Range("E99").Value = 307
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 31

# Row 109: Sudan del Sur - refreshed counters
$ws.Range("B109").Value = 1882
$ws.Range("C109").Value = 18
$ws.Range("E109").Value = 1726

# Row 122: Sierra Leona - refreshed counters
$ws.Range("B122").Value = 1309
$ws.Range("C122").Value = 11
$ws.Range("D122").Value = 746
$ws.Range("E122").Value = 510

# Row 127: Jordania - refreshed counters
$ws.Range("B127").Value = 1015
$ws.Range("C127").Value = 7
$ws.Range("D127").Value = 722
$ws.Range("E127").Value = 284

# Row 134: Republica del Chad - refreshed counters
$ws.Range("D134").Value = 746
$ws.Range("E134").Value = 38

# Row 202: Dominica -> Fiyi (rank swap)
$ws.Range("A202").Value = "Fiyi"

# Row 203: Fiyi -> Dominica (rank swap)
$ws.Range("A203").Value = "Dominica"

# Row 208: Islas Turcas y Caicos -> Santa Sede (rank swap)
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

# Row 209: Santa Sede -> Islas Turcas y Caicos (rank swap)
$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1

# Row 213: Islas Virgenes Britanicas -> Papua Nueva Guinea (rank swap)
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

# Row 214: Papua Nueva Guinea -> Islas Virgenes Britanicas (rank swap)
$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1

# Title cell: bump "datos actualizados" timestamp
$ws.Range("A1").Value = "Datos actualizados a 20 de Junio de 2020 a las 18:41"
